$d = $word.ActiveDocument

# Step 1: Expand the "5)" paragraph into the full sentence.
$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range
$r.Text = "5) Choose a solution and develop a plan to implement it."

# Step 2: Append the new paragraphs that follow it.
$newParas = @(
  "For this problem we are going to choose solution A. It provides a vast amount of flexibility, accuracy, and will be quick to use compared to other solutions.",
  "",
  "We’ll start with the first goal number of ten, as it’s easy to get to and check our work.",
  "1 = thumb. 2 = pointer. 3 = middle. 4 = ring. 5 = pinky.",
  "6 = ring. 7 = middle. 8 = pointer. 9 = thumb. 10 = pointer.",
  "",
  "The answer for 1-10 = pointer finger. Now, how do we get that into a formula? Well, right now we have an odd number as we go through each rotation, nine. This is caused by the pinky only being counted once, in our formula, we need to divide by an even number and we can do this by subtracting one thumb from the equation. This gives up a total of eight fingers. We divide each number by eight, indicating that each finger has been used an equal number of times. Afterwards, we look at the remainder of our division. Counting the thumb as one, we should simply be able to count up to the appropriate finger using the remainder to tell us which finger the number will land on.",
  "",
  "So 1-100 would be 100/8 = 12R4. So 100 will fall on the ring finger.",
  "For 1-1000, 1000/8 = 125R0. So 1000 will fall on the pinky finger."
)

foreach ($t in $newParas) {
  $cur = $d.Paragraphs.Last.Range
  $cur.InsertParagraphAfter()
  if ($t -ne "") {
    $newRange = $d.Paragraphs.Last.Range
    $newRange.Text = $t
  }
}
